$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from N1 to the new header cells O1:R1, then set their values
$ws.Range("N1").Copy()
$ws.Range("O1:R1").PasteSpecial(-4122)
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16

# Row 2
$ws.Range("C2").Value = 1.048680914502046
$ws.Range("D2").Value = 1.060465210219927
$ws.Range("E2").Value = 1.056685687889239
$ws.Range("F2").Value = 1.064601963877501
$ws.Range("I2").Value = 1.072986674555051
$ws.Range("J2").Value = 1.0694970641182
$ws.Range("K2").Value = 1.071141855701342
$ws.Range("L2").Value = 1.067407994866122
$ws.Range("M2").Value = 1.075229092811855
$ws.Range("O2").Value = 1.03
$ws.Range("P2").Value = 1.068110660655609
$ws.Range("Q2").Value = 1.02
$ws.Range("R2").Value = 1.061373196875256

# Row 3
$ws.Range("C3").Value = 1.053159196735994
$ws.Range("D3").Value = 1.063612145219397
$ws.Range("E3").Value = 1.060263458832078
$ws.Range("F3").Value = 1.067859489973878
$ws.Range("I3").Value = 1.074541994683033
$ws.Range("J3").Value = 1.072298602857183
$ws.Range("K3").Value = 1.073490995719344
$ws.Range("L3").Value = 1.070179227148596
$ws.Range("M3").Value = 1.077691924119687
$ws.Range("O3").Value = 1.03
$ws.Range("P3").Value = 1.070059753059539
$ws.Range("Q3").Value = 1.02
$ws.Range("R3").Value = 1.063031562611077

# Row 4
$ws.Range("C4").Value = 1.056002245051416
$ws.Range("D4").Value = 1.065613459674793
$ws.Range("E4").Value = 1.062539416767815
$ws.Range("F4").Value = 1.069933303952453
$ws.Range("I4").Value = 1.075520460033138
$ws.Range("J4").Value = 1.074074336936753
$ws.Range("K4").Value = 1.074979587481642
$ws.Range("L4").Value = 1.071937385573255
$ws.Range("M4").Value = 1.079255063613868
$ws.Range("O4").Value = 1.03
$ws.Range("P4").Value = 1.071296824085196
$ws.Range("Q4").Value = 1.02
$ws.Range("R4").Value = 1.064084930438177

# Row 5
$ws.Range("C5").Value = 1.057191578275498
$ws.Range("D5").Value = 1.06645386623083
$ws.Range("E5").Value = 1.063494272820452
$ws.Range("F5").Value = 1.070803204055735
$ws.Range("I5").Value = 1.075929753705122
$ws.Range("J5").Value = 1.074818868377506
$ws.Range("K5").Value = 1.075605169893842
$ws.Range("L5").Value = 1.072675413955265
$ws.Range("M5").Value = 1.079911023100314
$ws.Range("O5").Value = 1.03
$ws.Range("P5").Value = 1.071815950875584
$ws.Range("Q5").Value = 1.02
$ws.Range("R5").Value = 1.064534184815975

# Row 6
$ws.Range("C6").Value = 1.057398793289291
$ws.Range("D6").Value = 1.066603242755972
$ws.Range("E6").Value = 1.063662743986469
$ws.Range("F6").Value = 1.070956072630741
$ws.Range("I6").Value = 1.076003554240577
$ws.Range("J6").Value = 1.074951455943117
$ws.Range("K6").Value = 1.07571843500169
$ws.Range("L6").Value = 1.072807444818936
$ws.Range("M6").Value = 1.080027944547203
$ws.Range("O6").Value = 1.03
$ws.Range("P6").Value = 1.071908482554247
$ws.Range("Q6").Value = 1.02
$ws.Range("R6").Value = 1.064622732474383

# Row 7
$ws.Range("C7").Value = 1.056040555525961
$ws.Range("D7").Value = 1.065648356617821
$ws.Range("E7").Value = 1.062575761467243
$ws.Range("F7").Value = 1.069964647670999
$ws.Range("I7").Value = 1.075540721183964
$ws.Range("J7").Value = 1.074106194059803
$ws.Range("K7").Value = 1.075011366968057
$ws.Range("L7").Value = 1.071970584702685
$ws.Range("M7").Value = 1.079283344359717
$ws.Range("O7").Value = 1.03
$ws.Range("P7").Value = 1.071319205493886
$ws.Range("Q7").Value = 1.02
$ws.Range("R7").Value = 1.064126873938077

# Row 8
$ws.Range("C8").Value = 1.050234703889744
$ws.Range("D8").Value = 1.061566250746563
$ws.Range("E8").Value = 1.057933513839276
$ws.Range("F8").Value = 1.065735320873604
$ws.Range("I8").Value = 1.073537020527958
$ws.Range("J8").Value = 1.070479692807239
$ws.Range("K8").Value = 1.071972174584241
$ws.Range("L8").Value = 1.068382012177156
$ws.Range("M8").Value = 1.076092826532002
$ws.Range("O8").Value = 1.03
$ws.Range("P8").Value = 1.068794222766447
$ws.Range("Q8").Value = 1.02
$ws.Range("R8").Value = 1.061982481374202

# Row 9
$ws.Range("C9").Value = 1.039527192031288
$ws.Range("D9").Value = 1.054054552010771
$ws.Range("E9").Value = 1.04940079782746
$ws.Range("F9").Value = 1.057977623113123
$ws.Range("I9").Value = 1.069760730438799
$ws.Range("J9").Value = 1.06375740608601
$ws.Range("K9").Value = 1.066327545157862
$ws.Range("L9").Value = 1.061740185245992
$ws.Range("M9").Value = 1.070195219737496
$ws.Range("O9").Value = 1.03
$ws.Range("P9").Value = 1.064126822866088
$ws.Range("Q9").Value = 1.02
$ws.Range("R9").Value = 1.057988483997483

# Row 10
$ws.Range("C10").Value = 1.03218604827434
$ws.Range("D10").Value = 1.048955617363492
$ws.Range("E10").Value = 1.043623638648712
$ws.Range("F10").Value = 1.052752385462655
$ws.Range("I10").Value = 1.067160599252113
$ws.Range("J10").Value = 1.059178401013266
$ws.Range("K10").Value = 1.06249385244357
$ws.Range("L10").Value = 1.057247526630556
$ws.Range("M10").Value = 1.066230306479172
$ws.Range("O10").Value = 1.03
$ws.Range("P10").Value = 1.061039019951276
$ws.Range("Q10").Value = 1.02
$ws.Range("R10").Value = 1.055294314447683

# Row 11
$ws.Range("C11").Value = 1.029893205509096
$ws.Range("D11").Value = 1.047587892917697
$ws.Range("E11").Value = 1.042179210565309
$ws.Range("F11").Value = 1.051595590743557
$ws.Range("I11").Value = 1.066595642960516
$ws.Range("J11").Value = 1.058073181876934
$ws.Range("K11").Value = 1.061664907770488
$ws.Range("L11").Value = 1.05634723809126
$ws.Range("M11").Value = 1.065605944788176
$ws.Range("O11").Value = 1.03
$ws.Range("P11").Value = 1.060970094112264
$ws.Range("Q11").Value = 1.02
$ws.Range("R11").Value = 1.054740233989603

# Row 12
$ws.Range("C12").Value = 1.029377723398503
$ws.Range("D12").Value = 1.047389207203301
$ws.Range("E12").Value = 1.042038458871988
$ws.Range("F12").Value = 1.051576661736762
$ws.Range("I12").Value = 1.066594335727698
$ws.Range("J12").Value = 1.057988782739286
$ws.Range("K12").Value = 1.061661536514801
$ws.Range("L12").Value = 1.056402277383265
$ws.Range("M12").Value = 1.065778229741202
$ws.Range("O12").Value = 1.03
$ws.Range("P12").Value = 1.061425268218359
$ws.Range("Q12").Value = 1.02
$ws.Range("R12").Value = 1.05473785051693

# Row 13
$ws.Range("C13").Value = 1.030228443767409
$ws.Range("D13").Value = 1.048109782068669
$ws.Range("E13").Value = 1.04292483029588
$ws.Range("F13").Value = 1.052466895963355
$ws.Range("I13").Value = 1.067052232754982
$ws.Range("J13").Value = 1.058718225500848
$ws.Range("K13").Value = 1.06232874336579
$ws.Range("L13").Value = 1.057231976477636
$ws.Range("M13").Value = 1.066612617945093
$ws.Range("O13").Value = 1.03
$ws.Range("P13").Value = 1.06235602292817
$ws.Range("Q13").Value = 1.02
$ws.Range("R13").Value = 1.055207137589639

# Row 14
$ws.Range("C14").Value = 1.031440810686406
$ws.Range("D14").Value = 1.049029556039955
$ws.Range("E14").Value = 1.044007544567972
$ws.Range("F14").Value = 1.053500204395995
$ws.Range("I14").Value = 1.067579620024502
$ws.Range("J14").Value = 1.059595557379751
$ws.Range("K14").Value = 1.063098351972939
$ws.Range("L14").Value = 1.05816061545952
$ws.Range("M14").Value = 1.067494847669545
$ws.Range("O14").Value = 1.03
$ws.Range("P14").Value = 1.06322363825517
$ws.Range("Q14").Value = 1.02
$ws.Range("R14").Value = 1.055752636466802

# Row 15
$ws.Range("C15").Value = 1.032045821448942
$ws.Range("D15").Value = 1.049467411266301
$ws.Range("E15").Value = 1.044510296554582
$ws.Range("F15").Value = 1.053967342637157
$ws.Range("I15").Value = 1.067817279030087
$ws.Range("J15").Value = 1.060000254168469
$ws.Range("K15").Value = 1.063445993266461
$ws.Range("L15").Value = 1.058571405436913
$ws.Range("M15").Value = 1.067871873166375
$ws.Range("O15").Value = 1.03
$ws.Range("P15").Value = 1.063558566385761
$ws.Range("Q15").Value = 1.02
$ws.Range("R15").Value = 1.056004131024236

# Row 16
$ws.Range("C16").Value = 1.034996335278441
$ws.Range("D16").Value = 1.051495208026566
$ws.Range("E16").Value = 1.04678982714865
$ws.Range("F16").Value = 1.056032005849487
$ws.Range("I16").Value = 1.068849920567447
$ws.Range("J16").Value = 1.061814533400403
$ws.Range("K16").Value = 1.064958904231997
$ws.Range("L16").Value = 1.060328366824094
$ws.Range("M16").Value = 1.069424347740868
$ws.Range("O16").Value = 1.03
$ws.Range("P16").Value = 1.064747606806935
$ws.Range("Q16").Value = 1.02
$ws.Range("R16").Value = 1.057076834861739

# Row 17
$ws.Range("C17").Value = 1.036627695538662
$ws.Range("D17").Value = 1.052571901349982
$ws.Range("E17").Value = 1.047975234595334
$ws.Range("F17").Value = 1.057073706186534
$ws.Range("I17").Value = 1.069364372302884
$ws.Range("J17").Value = 1.062747614679878
$ws.Range("K17").Value = 1.065716934438543
$ws.Range("L17").Value = 1.061191309448873
$ws.Range("M17").Value = 1.070149926989035
$ws.Range("O17").Value = 1.03
$ws.Range("P17").Value = 1.065194463791157
$ws.Range("Q17").Value = 1.02
$ws.Range("R17").Value = 1.057615276705512

# Row 18
$ws.Range("C18").Value = 1.037241408730967
$ws.Range("D18").Value = 1.052889208990351
$ws.Range("E18").Value = 1.04828035211792
$ws.Range("F18").Value = 1.057277361941711
$ws.Range("I18").Value = 1.069454591171259
$ws.Range("J18").Value = 1.062967727600808
$ws.Range("K18").Value = 1.065853450528138
$ws.Range("L18").Value = 1.061314680354188
$ws.Range("M18").Value = 1.070175595561857
$ws.Range("O18").Value = 1.03
$ws.Range("P18").Value = 1.064982009923425
$ws.Range("Q18").Value = 1.02
$ws.Range("R18").Value = 1.05770048066517

# Row 19
$ws.Range("C19").Value = 1.036963835215304
$ws.Range("D19").Value = 1.052551115269807
$ws.Range("E19").Value = 1.047820251151803
$ws.Range("F19").Value = 1.056747023150921
$ws.Range("I19").Value = 1.069182577794713
$ws.Range("J19").Value = 1.062572717756904
$ws.Range("K19").Value = 1.065460455012324
$ws.Range("L19").Value = 1.060801221576656
$ws.Range("M19").Value = 1.069593514616759
$ws.Range("O19").Value = 1.03
$ws.Range("P19").Value = 1.064203919627367
$ws.Range("Q19").Value = 1.02
$ws.Range("R19").Value = 1.057428853548442

# Row 20
$ws.Range("C20").Value = 1.034157180301637
$ws.Range("D20").Value = 1.050345198615814
$ws.Range("E20").Value = 1.04519230993125
$ws.Range("F20").Value = 1.054165822800303
$ws.Range("I20").Value = 1.067881656835435
$ws.Range("J20").Value = 1.06043205927958
$ws.Range("K20").Value = 1.06355701227741
$ws.Range("L20").Value = 1.058484475438603
$ws.Range("M20").Value = 1.067318706558558
$ws.Range("O20").Value = 1.03
$ws.Range("P20").Value = 1.061890081292079
$ws.Range("Q20").Value = 1.02
$ws.Range("R20").Value = 1.056086941068611

# Row 21
$ws.Range("C21").Value = 1.028446683883425
$ws.Range("D21").Value = 1.046342850838728
$ws.Range("E21").Value = 1.040640772994546
$ws.Range("F21").Value = 1.050022886471619
$ws.Range("I21").Value = 1.065792755374225
$ws.Range("J21").Value = 1.056801395034098
$ws.Range("K21").Value = 1.06049639841149
$ws.Range("L21").Value = 1.054891038256087
$ws.Range("M21").Value = 1.064114732500397
$ws.Range("O21").Value = 1.03
$ws.Range("P21").Value = 1.059314661259406
$ws.Range("Q21").Value = 1.02
$ws.Range("R21").Value = 1.05392623694488

# Row 22
$ws.Range("C22").Value = 1.024808194357649
$ws.Range("D22").Value = 1.043802953755881
$ws.Range("E22").Value = 1.037761200975076
$ws.Range("F22").Value = 1.047416412478292
$ws.Range("I22").Value = 1.064460955104196
$ws.Range("J22").Value = 1.054496872509233
$ws.Range("K22").Value = 1.05855389693964
$ws.Range("L22").Value = 1.052619966806173
$ws.Range("M22").Value = 1.06210363930019
$ws.Range("O22").Value = 1.03
$ws.Range("P22").Value = 1.05772305434918
$ws.Range("Q22").Value = 1.02
$ws.Range("R22").Value = 1.052539704822501

# Row 23
$ws.Range("C23").Value = 1.026720991747712
$ws.Range("D23").Value = 1.045130006747932
$ws.Range("E23").Value = 1.039267630523759
$ws.Range("F23").Value = 1.048782122250652
$ws.Range("I23").Value = 1.065155397242282
$ws.Range("J23").Value = 1.055700754972557
$ws.Range("K23").Value = 1.059564115193379
$ws.Range("L23").Value = 1.053803603586651
$ws.Range("M23").Value = 1.063153515339825
$ws.Range("O23").Value = 1.03
$ws.Range("P23").Value = 1.058553941145305
$ws.Range("Q23").Value = 1.02
$ws.Range("R23").Value = 1.053244538609188

# Row 24
$ws.Range("C24").Value = 1.034123320595489
$ws.Range("D24").Value = 1.050289525686688
$ws.Range("E24").Value = 1.045122713844111
$ws.Range("F24").Value = 1.054090698808088
$ws.Range("I24").Value = 1.067835461769739
$ws.Range("J24").Value = 1.060368269842142
$ws.Range("K24").Value = 1.063487531177197
$ws.Range("L24").Value = 1.05840119812605
$ws.Range("M24").Value = 1.067230139594454
$ws.Range("O24").Value = 1.03
$ws.Range("P24").Value = 1.061780230946541
$ws.Range("Q24").Value = 1.02
$ws.Range("R24").Value = 1.056011135390976

# Row 25
$ws.Range("C25").Value = 1.042388696346644
$ws.Range("D25").Value = 1.056070999298257
$ws.Range("E25").Value = 1.05168651879342
$ws.Range("F25").Value = 1.060050699555417
$ws.Range("I25").Value = 1.070789079512773
$ws.Range("J25").Value = 1.065569129775487
$ws.Range("K25").Value = 1.06785705680247
$ws.Range("L25").Value = 1.06353215905917
$ws.Range("M25").Value = 1.071783211186503
$ws.Range("O25").Value = 1.03
$ws.Range("P25").Value = 1.065383571134507
$ws.Range("Q25").Value = 1.02
$ws.Range("R25").Value = 1.059097665391026
